$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Swap the Fecha (D) and Volumen (J) values between row 3 and row 5.
$d3 = $ws.Range("D3").Value2
$j3 = $ws.Range("J3").Value2
$d5 = $ws.Range("D5").Value2
$j5 = $ws.Range("J5").Value2

$ws.Range("D3").Value2 = $d5
$ws.Range("J3").Value2 = $j5
$ws.Range("D5").Value2 = $d3
$ws.Range("J5").Value2 = $j3
